$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Week 6 standings: row 7/8 swap places (Philip Milam overtakes Senay Semere),
# so update the names in column A for those two rows before writing new totals.
$ws.Range("A7").Value = "Philip Milam"
$ws.Range("A8").Value = "Senay Semere"

# Updated "Total Points" (column B) for week 6
$ws.Range("B2").Value = 1975.9
$ws.Range("B3").Value = 1721.5
$ws.Range("B4").Value = 1650.3
$ws.Range("B5").Value = 1387.6
$ws.Range("B6").Value = 1339.6
$ws.Range("B7").Value = 1174
$ws.Range("B8").Value = 909.1
$ws.Range("B9").Value = 743.8
$ws.Range("B10").Value = 621.5
$ws.Range("B11").Value = 182

# New empty helper cells in column H (rows 2-11) carrying a "0.0" number format
$ws.Range("H2:H11").NumberFormat = "0.0"

# Restore the active selection to D14
$ws.Range("D14").Select()
